$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: bold A1:AA1 (the real headers A1:D1 pick up a bold Arial
# variant, the filler cells E1:Y1 and the two new header cells Z1/AA1 pick
# up a bold Calibri variant - this matches the workbook's own default font) ---
$ws.Range("A1:AA1").Font.Bold = $true

# --- New columns: createdAt / lastUpdate ---
$ws.Range("Z1").Value = "createdAt"
$ws.Range("AA1").Value = "lastUpdate"

# --- Populate the (previously empty) id column A with the row index ---
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 5

# --- Timestamp formulas for the two new columns, one cell at a time so each
# keeps its own independent formula instead of being grouped into a shared
# formula block ---
for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 26).Formula = "=NOW()"
    $ws.Cells.Item($r, 27).Formula = "=NOW()"
}
$ws.Range("Z2:AA6").NumberFormat = "m/d/yy h:mm"

# --- Column widths for the new columns ---
$ws.Range("Z1").EntireColumn.ColumnWidth = 16.1
$ws.Range("AA1").EntireColumn.ColumnWidth = 16.1

# --- View: scroll so column I is the left-most visible column, and select
# the newly added lastUpdate column's data ---
$ws.Range("AA2:AA6").Select()
$excel.ActiveWindow.ScrollColumn = 9

Write-Output "done"
